$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format ("@") on the Price (D) column cells we are about to update.
# These cells hold plain text values (e.g. "16.70", "1.722.44"), and without this
# Excel would auto-convert numeric-looking strings to real numbers, stripping
# trailing zeros / grouping dots and changing the cell type away from text.
$priceCells = @(
    "D2",
    "D3",
    "D5",
    "D8",
    "D10",
    "D12",
    "D13",
    "D14",
    "D16",
    "D17",
    "D18",
    "D19",
    "D20",
    "D22",
    "D23",
    "D25",
    "D27",
    "D30",
    "D33",
    "D34",
    "D36",
    "D37",
    "D38",
    "D39",
    "D41",
    "D42",
    "D45",
    "D46",
    "D47",
    "D48",
    "D50"
)
foreach ($cellRef in $priceCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Apply the updated values cell by cell.
$ws.Range("D2").Value = "28.041.58"
$ws.Range("E2").Value = "  +3.42%  "
$ws.Range("D3").Value = "1.722.44"
$ws.Range("E3").Value = "  +2.51%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "218.85"
$ws.Range("E5").Value = "  +1.73%  "
$ws.Range("E6").Value = "  +0.59%  "
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").Value = "24.29"
$ws.Range("E8").Value = "  +14.26%  "
$ws.Range("E9").Value = "  +2.95%  "
$ws.Range("D10").Value = "0.0633"
$ws.Range("E10").Value = "  +1.62%  "
$ws.Range("E11").Value = "  +1.56%  "
$ws.Range("D12").Value = "1.965.43"
$ws.Range("E12").Value = "  +2.50%  "
$ws.Range("D13").Value = "1.723.67"
$ws.Range("E13").Value = "  +2.48%  "
$ws.Range("D14").Value = "4.28"
$ws.Range("E14").Value = "  +3.13%  "
$ws.Range("E15").Value = "  +4.67%  "
$ws.Range("D16").Value = "67.53"
$ws.Range("E16").Value = "  +2.05%  "
$ws.Range("D17").Value = "27.996.70"
$ws.Range("E17").Value = "  +3.26%  "
$ws.Range("D18").Value = "242.27"
$ws.Range("E18").Value = "  +1.83%  "
$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").Value = "8.00"
$ws.Range("E19").Value = "  -1.39%  "
$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").Value = "0.0₃0754"
$ws.Range("E20").Value = "  +1.23%  "
$ws.Range("E21").Value = "  +0.01%  "
$ws.Range("D22").Value = "4.63"
$ws.Range("E22").Value = "  +2.50%  "
$ws.Range("D23").Value = "9.66"
$ws.Range("E23").Value = "  +2.14%  "
$ws.Range("E24").Value = "  +0.26%  "
$ws.Range("D25").Value = "148.66"
$ws.Range("E26").Value = "  +3.51%  "
$ws.Range("D27").Value = "16.70"
$ws.Range("E27").Value = "  +2.41%  "
$ws.Range("E28").Value = "  +1.03%  "
$ws.Range("E29").Value = "  -0.34%  "
$ws.Range("D30").Value = "0.0509"
$ws.Range("E30").Value = "  +1.89%  "
$ws.Range("E31").Value = "  +1.78%  "
$ws.Range("E32").Value = "  +2.11%  "
$ws.Range("B33").Value = "Maker"
$ws.Range("C33").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D33").Value = "1.490.94"
$ws.Range("E33").Value = "  -4.45%  "
$ws.Range("B34").Value = "InternetComputer(DFINITY)"
$ws.Range("C34").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D34").Value = "3.26"
$ws.Range("E34").Value = "  +2.09%  "
$ws.Range("E35").Value = "  -2.51%  "
$ws.Range("D36").Value = "0.955"
$ws.Range("E36").Value = "  +2.22%  "
$ws.Range("D37").Value = "0.606"
$ws.Range("E37").Value = "  +0.38%  "
$ws.Range("D38").Value = "2.41"
$ws.Range("E38").Value = "  +1.04%  "
$ws.Range("D39").Value = "0.0175"
$ws.Range("E39").Value = "  +0.11%  "
$ws.Range("E40").Value = "  +1.59%  "
$ws.Range("D41").Value = "70.55"
$ws.Range("E41").Value = "  +2.29%  "
$ws.Range("D42").Value = "5.84"
$ws.Range("E42").Value = "  +2.99%  "
$ws.Range("E43").Value = "  -0.01%  "
$ws.Range("E44").Value = "  +1.85%  "
$ws.Range("D45").Value = "1.869.59"
$ws.Range("E45").Value = "  +2.35%  "
$ws.Range("D46").Value = "0.800"
$ws.Range("E46").Value = "  +2.03%  "
$ws.Range("D47").Value = "1.77"
$ws.Range("E47").Value = "  +11.93%  "
$ws.Range("D48").Value = "90.99"
$ws.Range("E48").Value = "  +0.28%  "
$ws.Range("E49").Value = "  +4.29%  "
$ws.Range("D50").Value = "8.26"
$ws.Range("E50").Value = "  +2.06%  "
$ws.Range("E51").Value = "  +0.65%  "
